$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume data updates (refresh run)
$ws.Range('D2').Value = '27.393.05'
$ws.Range('E2').Value = '  -3.41%  '
$ws.Range('D3').Value = '1.650.66'
$ws.Range('E3').Value = '  -3.59%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.99'
$ws.Range('E5').Value = '  -1.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.508'
$ws.Range('E6').Value = '  -2.75%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.75%  '
$ws.Range('E9').Value = '  -1.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0616'
$ws.Range('E10').Value = '  -2.62%  '
$ws.Range('E11').Value = '  -1.51%  '
$ws.Range('D12').Value = '1.884.93'
$ws.Range('E12').Value = '  -3.63%  '
$ws.Range('D13').Value = '1.650.17'
$ws.Range('E13').Value = '  -3.70%  '
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.564'
$ws.Range('E15').Value = '  +0.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.72'
$ws.Range('E16').Value = '  -2.48%  '
$ws.Range('D17').Value = '27.399.62'
$ws.Range('E17').Value = '  -3.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '234.59'
$ws.Range('E18').Value = '  -7.84%  '
$ws.Range('E19').Value = '  -2.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.50'
$ws.Range('E20').Value = '  -3.18%  '
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('E22').Value = '  -3.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.28'
$ws.Range('E23').Value = '  -3.26%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.07'
$ws.Range('E25').Value = '  -1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.18'
$ws.Range('E26').Value = '  -2.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.06'
$ws.Range('E27').Value = '  -3.53%  '
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('E29').Value = '  -2.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0497'
$ws.Range('E30').Value = '  -2.77%  '
$ws.Range('E31').Value = '  -0.85%  '
$ws.Range('E32').Value = '  -2.84%  '
$ws.Range('D33').Value = '1.465.60'
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('E34').Value = '  -3.80%  '
$ws.Range('E35').Value = '  -4.95%  '
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.911'
$ws.Range('E37').Value = '  -5.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.571'
$ws.Range('E38').Value = '  -4.37%  '
$ws.Range('E39').Value = '  -2.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.04'
$ws.Range('E40').Value = '  -0.32%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.44'
$ws.Range('E42').Value = '  -3.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.23'
$ws.Range('E43').Value = '  -6.24%  '
$ws.Range('E44').Value = '  -2.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.786'
$ws.Range('E45').Value = '  -1.83%  '
$ws.Range('D46').Value = '1.793.05'
$ws.Range('E46').Value = '  -3.60%  '
$ws.Range('E47').Value = '  -0.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.18'
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('E49').Value = '  -4.70%  '
$ws.Range('E50').Value = '  -2.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.78'
$ws.Range('E51').Value = '  -3.52%  '
